$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Client")

# Update the client name values (2607 -> 2707)
$ws.Range("B2").Value = "Anh Tester Client 2707A1"
$ws.Range("B3").Value = "Anh Tester Client 2707A2"
$ws.Range("B4").Value = "Anh Tester Client 2707A3"

# Move the active selection from B7 to B6
$ws.Activate()
$ws.Range("B6").Select()
